# Updated RAD Notice Number Error Message.
# Refresh the "Date" column (column B) timestamps for the test rows that
# were re-run, leaving the DONOTRUN rows (19-24) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @{
    2  = "Wed Dec 20 12:43:08 EST 2023"
    3  = "Wed Dec 20 12:43:18 EST 2023"
    4  = "Wed Dec 20 12:43:28 EST 2023"
    5  = "Wed Dec 20 12:43:38 EST 2023"
    6  = "Wed Dec 20 12:43:48 EST 2023"
    7  = "Wed Dec 20 12:43:58 EST 2023"
    8  = "Wed Dec 20 12:44:08 EST 2023"
    9  = "Wed Dec 20 12:44:18 EST 2023"
    10 = "Wed Dec 20 12:44:28 EST 2023"
    11 = "Wed Dec 20 12:44:38 EST 2023"
    12 = "Wed Dec 20 12:44:48 EST 2023"
    13 = "Wed Dec 20 12:44:58 EST 2023"
    14 = "Wed Dec 20 12:45:09 EST 2023"
    15 = "Wed Dec 20 12:45:19 EST 2023"
    16 = "Wed Dec 20 12:45:29 EST 2023"
    17 = "Wed Dec 20 12:45:39 EST 2023"
    18 = "Wed Dec 20 12:45:49 EST 2023"
    25 = "Wed Dec 20 12:46:00 EST 2023"
    26 = "Wed Dec 20 12:46:09 EST 2023"
    27 = "Wed Dec 20 12:46:19 EST 2023"
    28 = "Wed Dec 20 12:46:29 EST 2023"
    29 = "Wed Dec 20 12:46:39 EST 2023"
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 2).Value = $dates[$row]
}
